# Apply the "Update countries & provincias Spain" data refresh.
#
# The source COVID dashboard re-ran its data pull; several countries' case
# counts moved enough to change their sort position (table is kept sorted
# descending by "Casos totales" / column B). Angola and Aruba both grew past
# their neighbours, so those rows move up while the other affected rows keep
# their own original data, just shifted down one position. All other affected
# rows are simple in-place numeric refreshes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp (row 1, col A)
$ws.Range("A1").Value = "Datos actualizados a 20 de Agosto de 2020 a las 22:12"

# Row 4
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 5729773
$ws.Cells.Item(4, 3).Value = 28842
$ws.Cells.Item(4, 4).Value = 3072781
$ws.Cells.Item(4, 5).Value = 2480069
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 589
$ws.Cells.Item(4, 8).Value = 176923

# Row 12
$ws.Cells.Item(12, 1).Value = "Chile"
$ws.Cells.Item(12, 2).Value = 391849
$ws.Cells.Item(12, 3).Value = 1812
$ws.Cells.Item(12, 4).Value = 366063
$ws.Cells.Item(12, 5).Value = 15115
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value = 93
$ws.Cells.Item(12, 8).Value = 10671

# Row 22
$ws.Cells.Item(22, 1).Value = "Alemania"
$ws.Cells.Item(22, 2).Value = 231188
$ws.Cells.Item(22, 3).Value = 1488
$ws.Cells.Item(22, 4).Value = 204800
$ws.Cells.Item(22, 5).Value = 17064
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 10
$ws.Cells.Item(22, 8).Value = 9324

# Row 59
$ws.Cells.Item(59, 1).Value = "Suiza"
$ws.Cells.Item(59, 2).Value = 39026
$ws.Cells.Item(59, 3).Value = 266
$ws.Cells.Item(59, 4).Value = 33900
$ws.Cells.Item(59, 5).Value = 3128
$ws.Cells.Item(59, 6).Value = 0
$ws.Cells.Item(59, 7).Value = 2
$ws.Cells.Item(59, 8).Value = 1998

# Row 67
$ws.Cells.Item(67, 1).Value = "Costa Rica"
$ws.Cells.Item(67, 2).Value = 31075
$ws.Cells.Item(67, 3).Value = 666
$ws.Cells.Item(67, 4).Value = 9939
$ws.Cells.Item(67, 5).Value = 20803
$ws.Cells.Item(67, 6).Value = 0
$ws.Cells.Item(67, 7).Value = 12
$ws.Cells.Item(67, 8).Value = 333

# Row 103
$ws.Cells.Item(103, 1).Value = "Mauritania"
$ws.Cells.Item(103, 2).Value = 6848
$ws.Cells.Item(103, 3).Value = 19
$ws.Cells.Item(103, 4).Value = 6123
$ws.Cells.Item(103, 5).Value = 567
$ws.Cells.Item(103, 6).Value = 0
$ws.Cells.Item(103, 7).Value = 0
$ws.Cells.Item(103, 8).Value = 158

# Row 105
$ws.Cells.Item(105, 1).Value = "Zimbabue"
$ws.Cells.Item(105, 2).Value = 5745
$ws.Cells.Item(105, 3).Value = 102
$ws.Cells.Item(105, 4).Value = 4525
$ws.Cells.Item(105, 5).Value = 1069
$ws.Cells.Item(105, 6).Value = 0
$ws.Cells.Item(105, 7).Value = 1
$ws.Cells.Item(105, 8).Value = 151

# Row 107
$ws.Cells.Item(107, 1).Value = "Malaui"
$ws.Cells.Item(107, 2).Value = 5282
$ws.Cells.Item(107, 3).Value = 42
$ws.Cells.Item(107, 4).Value = 2883
$ws.Cells.Item(107, 5).Value = 2234
$ws.Cells.Item(107, 6).Value = 0
$ws.Cells.Item(107, 7).Value = 1
$ws.Cells.Item(107, 8).Value = 165

# Row 119
$ws.Cells.Item(119, 1).Value = "Cabo Verde"
$ws.Cells.Item(119, 2).Value = 3368
$ws.Cells.Item(119, 3).Value = 47
$ws.Cells.Item(119, 4).Value = 2462
$ws.Cells.Item(119, 5).Value = 869
$ws.Cells.Item(119, 6).Value = 0
$ws.Cells.Item(119, 7).Value = 1
$ws.Cells.Item(119, 8).Value = 37

# Row 136
$ws.Cells.Item(136, 1).Value = "Angola"
$ws.Cells.Item(136, 2).Value = 2044
$ws.Cells.Item(136, 3).Value = 29
$ws.Cells.Item(136, 4).Value = 742
$ws.Cells.Item(136, 5).Value = 1209
$ws.Cells.Item(136, 6).Value = 0
$ws.Cells.Item(136, 7).Value = 1
$ws.Cells.Item(136, 8).Value = 93

# Row 137
$ws.Cells.Item(137, 1).Value = "Islandia"
$ws.Cells.Item(137, 2).Value = 2040
$ws.Cells.Item(137, 3).Value = 5
$ws.Cells.Item(137, 4).Value = 1913
$ws.Cells.Item(137, 5).Value = 117
$ws.Cells.Item(137, 6).Value = 0
$ws.Cells.Item(137, 7).Value = 0
$ws.Cells.Item(137, 8).Value = 10

# Row 147
$ws.Cells.Item(147, 1).Value = "Aruba"
$ws.Cells.Item(147, 2).Value = 1387
$ws.Cells.Item(147, 3).Value = 91
$ws.Cells.Item(147, 4).Value = 277
$ws.Cells.Item(147, 5).Value = 1104
$ws.Cells.Item(147, 6).Value = 0
$ws.Cells.Item(147, 7).Value = 1
$ws.Cells.Item(147, 8).Value = 6

# Row 148
$ws.Cells.Item(148, 1).Value = "Republica de Chipre"
$ws.Cells.Item(148, 2).Value = 1385
$ws.Cells.Item(148, 3).Value = 0
$ws.Cells.Item(148, 4).Value = 878
$ws.Cells.Item(148, 5).Value = 487
$ws.Cells.Item(148, 6).Value = 0
$ws.Cells.Item(148, 7).Value = 0
$ws.Cells.Item(148, 8).Value = 20

# Row 149
$ws.Cells.Item(149, 1).Value = "Georgia"
$ws.Cells.Item(149, 2).Value = 1370
$ws.Cells.Item(149, 3).Value = 9
$ws.Cells.Item(149, 4).Value = 1108
$ws.Cells.Item(149, 5).Value = 245
$ws.Cells.Item(149, 6).Value = 0
$ws.Cells.Item(149, 7).Value = 0
$ws.Cells.Item(149, 8).Value = 17

# Row 150
$ws.Cells.Item(150, 1).Value = "Letonia"
$ws.Cells.Item(150, 2).Value = 1327
$ws.Cells.Item(150, 3).Value = 1
$ws.Cells.Item(150, 4).Value = 1093
$ws.Cells.Item(150, 5).Value = 201
$ws.Cells.Item(150, 6).Value = 0
$ws.Cells.Item(150, 7).Value = 0
$ws.Cells.Item(150, 8).Value = 33

# Row 151
$ws.Cells.Item(151, 1).Value = "Botsuana"
$ws.Cells.Item(151, 2).Value = 1308
$ws.Cells.Item(151, 3).Value = 0
$ws.Cells.Item(151, 4).Value = 136
$ws.Cells.Item(151, 5).Value = 1169
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 7).Value = 0
$ws.Cells.Item(151, 8).Value = 3

# Row 167
$ws.Cells.Item(167, 1).Value = "Belice"
$ws.Cells.Item(167, 2).Value = 605
$ws.Cells.Item(167, 3).Value = 52
$ws.Cells.Item(167, 4).Value = 38
$ws.Cells.Item(167, 5).Value = 562
$ws.Cells.Item(167, 6).Value = 0
$ws.Cells.Item(167, 7).Value = 1
$ws.Cells.Item(167, 8).Value = 5

